# New Submission Synced: 2026-02-08 17:19:49
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

# C3 was stored as a text value ("22"); normalize it to a real number to
# match the rest of the numeric columns.
$ws.Range("C3").Value = 22

# Append the new submission row.
$ws.Range("A4").Value = "2026-02-08 17:19:49"
$ws.Range("B4").Value = "Ummu Kalthum Muhammad Yahaya"
# Admission No for this submission keeps its text representation (matches
# the sheet's historical mixed-typed column), so force it in as text.
$ws.Range("C4").Value = "'43"
$ws.Range("D4").Value = 10
